# Updates "Horarios 141" workbook with the latest scrape results.
# New scrape timestamp: 04:40:48 (was 04:18:52)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "LP1912": refreshed header + 7 new rows (14 -> 21 total rows)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:40:48"
$ws1.Range("A3").Value = "Total filas: 21"

$ws1.Range("A17").Value = "04:40:48"
$ws1.Range("B17").Value = "05:54"
$ws1.Range("C17").Value = "10_OLMOS"
$ws1.Range("D17").Value = 74
$ws1.Range("E17").Value = "LP1912"

$ws1.Range("A18").Value = "04:40:48"
$ws1.Range("B18").Value = "06:04"
$ws1.Range("C18").Value = "16_SANTA ANA"
$ws1.Range("D18").Value = 84
$ws1.Range("E18").Value = "LP1912"

$ws1.Range("A19").Value = "04:40:48"
$ws1.Range("B19").Value = "06:05"
$ws1.Range("C19").Value = "16_SANTA ANA"
$ws1.Range("D19").Value = 107
$ws1.Range("E19").Value = "LP1912"

$ws1.Range("A20").Value = "04:18:52"
$ws1.Range("B20").Value = "06:11"
$ws1.Range("C20").Value = "215A_EL PATO"
$ws1.Range("D20").Value = 113
$ws1.Range("E20").Value = "LP1912"

$ws1.Range("A21").Value = "04:18:52"
$ws1.Range("B21").Value = "06:13"
$ws1.Range("C21").Value = "225_HARAS DEL SUR"
$ws1.Range("D21").Value = 115
$ws1.Range("E21").Value = "LP1912"

$ws1.Range("A22").Value = "04:40:48"
$ws1.Range("B22").Value = "06:14"
$ws1.Range("C22").Value = "225_HARAS DEL SUR"
$ws1.Range("D22").Value = 94
$ws1.Range("E22").Value = "LP1912"

$ws1.Range("A23").Value = "04:40:48"
$ws1.Range("B23").Value = "06:21"
$ws1.Range("C23").Value = "26_HERNANDEZ"
$ws1.Range("D23").Value = 101
$ws1.Range("E23").Value = "LP1912"

$ws1.Range("A24").Value = "04:40:48"
$ws1.Range("B24").Value = "06:27"
$ws1.Range("C24").Value = "23_HERNANDEZ"
$ws1.Range("D24").Value = 107
$ws1.Range("E24").Value = "LP1912"

$ws1.Range("A25").Value = "04:40:48"
$ws1.Range("B25").Value = "06:29"
$ws1.Range("C25").Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Range("D25").Value = 109
$ws1.Range("E25").Value = "LP1912"

$ws1.Range("A26").Value = "04:40:48"
$ws1.Range("B26").Value = "06:31"
$ws1.Range("C26").Value = "16_SANTA ANA"
$ws1.Range("D26").Value = 111
$ws1.Range("E26").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "LP1912-215": only the refreshed timestamp changes
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 04:40:48"

# ---------------------------------------------------------------------
# Sheet "6203-6173": refreshed header + 1 new row (3 -> 4 total rows)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:40:48"
$ws3.Range("A3").Value = "Total filas: 4"

$ws3.Range("A9").Value = "04:40:48"
$ws3.Range("B9").Value = "06:32"
$ws3.Range("C9").Value = "215C_LA PLATA"
$ws3.Range("D9").Value = 112
$ws3.Range("E9").Value = "L6203"
